$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 2.55
$ws.Range("L2").Value = 3.4
$ws.Range("Z2").Value = 34
$ws.Range("AG2").Value = 6
$ws.Range("AW2").Value = 4.33
$ws.Range("AX2").Value = 15

# Row 3 updates
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 7.5
$ws.Range("AN3").Value = 3.75
$ws.Range("AX3").Value = 23
$ws.Range("BB3").Value = 301
